$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.133226633071899
$ws.Range("B1").Value = 2.649789333343506
$ws.Range("C1").Value = 3.942341327667236
$ws.Range("D1").Value = 3.685919046401978
$ws.Range("E1").Value = 1.234290599822998
